# Update the cryptocurrency price/volume list.
# Periodic GitHub Actions refresh of the "cryptos" sheet: most rows keep the
# same coin/link but get new Price / Volume(1h) values. Rows 13 and 14 swap
# which coin (Polkadot / WrappedEther) occupies which row, each getting its
# own updated price/volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.343.59"
$ws.Range("E2").Value = "  +2.78%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.869.27"
$ws.Range("E3").Value = "  +1.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.12%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.76"
$ws.Range("E5").Value = "  +1.99%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.13%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4704"
$ws.Range("E7").Value = "  +1.88%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3957"
$ws.Range("E8").Value = "  +2.86%  "

# Row 9 - OKB
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.40"
$ws.Range("E9").Value = "  +3.51%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07992"
$ws.Range("E10").Value = "  +1.07%  "

# Row 11 - Polygon
$ws.Range("E11").Value = "  +1.63%  "

# Row 12 - Solana
$ws.Range("E12").Value = "  +2.34%  "

# Row 13 - WrappedEther (was Polkadot)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.872.02"
$ws.Range("E13").Value = "  +1.11%  "

# Row 14 - Polkadot (was WrappedEther)
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.009"
$ws.Range("E14").Value = "  +1.88%  "

# Row 15 - Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.282"
$ws.Range("E15").Value = "  +2.79%  "

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.18"
$ws.Range("E16").Value = "  +2.97%  "

# Row 17 - BinanceUSD
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.19%  "

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001038"
$ws.Range("E18").Value = "  +0.48%  "

# Row 19 - TRON
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06605"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20 - Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.70"
$ws.Range("E20").Value = "  +4.33%  "

# Row 21 - Dai
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.12%  "

# Row 22 - WrappedBTC
$ws.Range("D22").Value = "28.362.34"
$ws.Range("E22").Value = "  +2.82%  "

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.463"
$ws.Range("E23").Value = "  +1.76%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  +1.45%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.293"
$ws.Range("E25").Value = "  -0.52%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").Value = "2.097.77"
$ws.Range("E26").Value = "  +1.33%  "

# Row 27 - Monero
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.05"
$ws.Range("E27").Value = "  +1.43%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +2.09%  "

# Row 29 - LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.152"
$ws.Range("E29").Value = "  +3.24%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.494"
$ws.Range("E30").Value = "  +1.97%  "

# Row 31 - BitcoinCash
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.32"
$ws.Range("E31").Value = "  +0.62%  "

# Row 32 - ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9751"
$ws.Range("E32").Value = "  +0.39%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  +1.58%  "

# Row 34 - HuobiToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.591"
$ws.Range("E34").Value = "  +0.37%  "

# Row 35 - ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.384"
$ws.Range("E35").Value = "  +3.17%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  +1.71%  "

# Row 37 - VeChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02275"
$ws.Range("E37").Value = "  +2.68%  "

# Row 38 - Hedera
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06098"
$ws.Range("E38").Value = "  +1.67%  "

# Row 39 - FraxShare
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.475"
$ws.Range("E39").Value = "  +2.35%  "

# Row 40 - TrustWalletToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.182"
$ws.Range("E40").Value = "  +0.40%  "

# Row 41 - TheSandbox
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5985"
$ws.Range("E41").Value = "  +1.79%  "

# Row 42 - Frax
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.08%  "

# Row 43 - Algorand
$ws.Range("E43").Value = "  +1.42%  "

# Row 44 - Aptos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.40"
$ws.Range("E44").Value = "  +1.30%  "

# Row 45 - WEMIXTOKEN
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.305"
$ws.Range("E45").Value = "  +3.96%  "

# Row 46 - Decentraland
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5626"
$ws.Range("E46").Value = "  +0.99%  "

# Row 47 - EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.15"
$ws.Range("E47").Value = "  +0.68%  "

# Row 48 - NEARProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.969"
$ws.Range("E48").Value = "  +4.08%  "

# Row 49 - Cronos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06926"
$ws.Range("E49").Value = "  +3.55%  "

# Row 50 - Quant
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.30"
$ws.Range("E50").Value = "  +0.86%  "

# Row 51 - RenderToken
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.025"
$ws.Range("E51").Value = "  +13.74%  "
